$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert a brand-new "Instrument" column before the old column I (Sector *),
# pushing every column from I onward one slot to the right (I->J, J->K, ...).
# ---------------------------------------------------------------------------
$ws.Columns("I:I").Insert()

# New column inherits the look of its left neighbour (old column H) as Excel
# normally does when inserting a column - restore that width as closely as
# the host lets us.
$ws.Columns("I:I").ColumnWidth = 7.33

# Header + body values for the freshly inserted column.
$ws.Range("I1").Value = "Instrument"
$ws.Range("I2:I9").Value = "Stock"

# Give the new body cells the same cell style the rest of the data rows use
# (copy format only, from the neighbouring data column).
$ws.Range("J2").Copy()
$ws.Range("I2:I9").PasteSpecial(-4122)
$dummy = $excel.CutCopyMode
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# The header comments that used to sit on I1..N1 need to follow their cells
# to J1..O1. Walk right-to-left so we never overwrite a comment we still
# need to read.
# ---------------------------------------------------------------------------
$oldCols = @("N", "M", "L", "K", "J", "I")
$newCols = @("O", "N", "M", "L", "K", "J")
for ($i = 0; $i -lt $oldCols.Length; $i++) {
    $oldRef = $oldCols[$i] + "1"
    $newRef = $newCols[$i] + "1"
    $cmt = $ws.Range($oldRef).Comment
    if ($cmt -ne $null) {
        $txt = $cmt.Text()
        $cmt.Delete()
        $ws.Range($newRef).AddComment($txt)
    }
}

# ---------------------------------------------------------------------------
# Re-create the three list validations so their sqref collapses back down to
# a single contiguous range (matching how Excel re-writes them after the
# column shift) instead of staying split as "…2:…9 …10:…1048576".
# ---------------------------------------------------------------------------
$gValidation = $ws.Range("G2:G1048576")
$gValidation.Validation.Delete()
$gValidation.Validation.Add(3, 1, 1, '"Pool,CoInvest"')

$mValidation = $ws.Range("M2:M1048576")
$mValidation.Validation.Delete()
$mValidation.Validation.Add(3, 1, 1, '"Yes,No"')

$nValidation = $ws.Range("N2:N1048576")
$nValidation.Validation.Delete()
$nValidation.Validation.Add(3, 1, 1, '"Domestic,Overseas"')

# ---------------------------------------------------------------------------
# Selection moves onto the new Instrument column's data cells.
# ---------------------------------------------------------------------------
$dummySelect = $ws.Range("I3:I9").Select()
